# Apply the "most recent xls files" update:
#  - rename the sheet (and its matching defined name) from
#    Blood_Lead_Level -> Blood_Lead
#  - rename the "...5yavg" header labels to "...c1115" in the 5 summary
#    columns (G1, M1, S1, Y1, AE1)
#  - narrow column A slightly
#  - set the worksheet to print in portrait orientation

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- sheet + defined name rename -----------------------------------------
$ws.Name = "Blood_Lead"
$wb.Names.Item("Blood_Lead_Level").Name = "Blood_Lead"

# --- header label text updates --------------------------------------------
$ws.Range("G1").Value  = "_ebll_c1115"
$ws.Range("M1").Value  = "_w_ebll_c1115"
$ws.Range("S1").Value  = "_b_ebll_c1115"
$ws.Range("Y1").Value  = "_a_ebll_c1115"
$ws.Range("AE1").Value = "_o_ebll_c1115"

# --- column A width ---------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 26.5

# --- page setup: portrait orientation --------------------------------
$ws.PageSetup.Orientation = 1
